# "Generate Report for Handback"
#
# The localization-status report is regenerated after a handback completes:
#   - Row 2 / Row 3 status text moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" on both the zh-cn and de-de sheets.
#   - The "Latest Target File" (E) and "Latest Handback File" (F) columns -
#     previously empty - are now populated with the handed-back file names,
#     each as a hyperlink (mirroring the existing A / C hyperlink cells).
#   - The "Latest Handback DateTime" (G) column moves from the placeholder
#     "0001-01-01 00:00:00" to the real handback timestamp.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$sheets = @(
    @{
        Name = "zh-cn"
        XlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eb3753161b6db667f0fda971f712a0fb1b3c2a40/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/026cd78d-c75f-434c-9286-2379c3ec760f.06e138f3c8b4177e1abca4892cfa570bc49e53a9.zh-cn.xlf"
        XlfName   = "026cd78d-c75f-434c-9286-2379c3ec760f.06e138f3c8b4177e1abca4892cfa570bc49e53a9.zh-cn.xlf"
        HandbackDateTimeRow2 = "2016-03-08 06:34:35"
        HandbackDateTimeRow3 = "2016-03-08 06:34:35"
    },
    @{
        Name = "de-de"
        XlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7de7f646fccdc52bc10498b6561e63b237d6f1d5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/026cd78d-c75f-434c-9286-2379c3ec760f.06e138f3c8b4177e1abca4892cfa570bc49e53a9.de-de.xlf"
        XlfName   = "026cd78d-c75f-434c-9286-2379c3ec760f.06e138f3c8b4177e1abca4892cfa570bc49e53a9.de-de.xlf"
        HandbackDateTimeRow2 = "2016-03-08 06:34:41"
        HandbackDateTimeRow3 = "2016-03-08 06:34:41"
    }
)

$mdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/654bafe25c73634c5bcd3c42f613a029e9c40485/e2e/026cd78d-c75f-434c-9286-2379c3ec760f.md"
$mdName = "026cd78d-c75f-434c-9286-2379c3ec760f.md"

foreach ($info in $sheets) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Row 2 + Row 3: status text updated now that handback is complete.
    $ws.Range("B2").Value = $newStatus
    $ws.Range("B3").Value = $newStatus

    # Row 2 + Row 3: "Latest Target File" (E) and "Latest Handback File" (F)
    # now carry the handed-back file references, as hyperlinks just like the
    # existing Source File Name (A) / Latest Handoff File (C) columns.
    $ws.Hyperlinks.Add($ws.Range("E2"), $mdTarget, "", "", $mdName)
    $ws.Hyperlinks.Add($ws.Range("F2"), $info.XlfTarget, "", "", $info.XlfName)
    $ws.Hyperlinks.Add($ws.Range("E3"), $mdTarget, "", "", $mdName)
    $ws.Hyperlinks.Add($ws.Range("F3"), $info.XlfTarget, "", "", $info.XlfName)

    # Row 2 + Row 3: "Latest Handback DateTime" (G) now holds the real
    # handback timestamp instead of the 0001-01-01 00:00:00 placeholder.
    $ws.Range("G2").Value = $info.HandbackDateTimeRow2
    $ws.Range("G3").Value = $info.HandbackDateTimeRow3
}
